# Updates cryptos list prices/volumes (GitHub Actions scheduled refresh).
# Numeric-looking "Price" strings are entered with a leading apostrophe so
# Excel keeps them as text (matching the source file's inline-string cells)
# instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.733.22"
$ws.Range("E2").Value = "  -1.38%  "
$ws.Range("D3").Value = "2.210.85"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").Value = "'240.33"
$ws.Range("E5").Value = "  -1.76%  "
$ws.Range("D6").Value = "'0.620"
$ws.Range("E6").Value = "  -1.42%  "
$ws.Range("D7").Value = "'72.22"
$ws.Range("E7").Value = "  -1.86%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "'0.586"
$ws.Range("E9").Value = "  -4.50%  "
$ws.Range("D10").Value = "'41.13"
$ws.Range("E10").Value = "  -3.57%  "
$ws.Range("D11").Value = "'0.0936"
$ws.Range("E11").Value = "  -2.97%  "
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("D13").Value = "'6.82"
$ws.Range("E13").Value = "  -4.27%  "
$ws.Range("D14").Value = "2.546.28"
$ws.Range("E14").Value = "  -0.36%  "
$ws.Range("D15").Value = "'14.07"
$ws.Range("E15").Value = "  -2.23%  "
$ws.Range("D16").Value = "'0.822"
$ws.Range("E16").Value = "  -3.27%  "
$ws.Range("D17").Value = "2.226.91"
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").Value = "41.637.49"
$ws.Range("E18").Value = "  -1.11%  "
$ws.Range("D19").Value = "'0.0000103"
$ws.Range("E19").Value = "  -7.46%  "
$ws.Range("D20").Value = "'6.10"
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("D21").Value = "'71.37"
$ws.Range("E21").Value = "  -0.98%  "
$ws.Range("D22").Value = "'10.66"
$ws.Range("E22").Value = "  +5.86%  "
$ws.Range("D23").Value = "'227.48"
$ws.Range("E23").Value = "  -1.57%  "
$ws.Range("D24").Value = "'2.01"
$ws.Range("E24").Value = "  -6.70%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").Value = "'11.23"
$ws.Range("E26").Value = "  -5.45%  "
$ws.Range("E27").Value = "  -0.46%  "
$ws.Range("E28").Value = "  -2.11%  "
$ws.Range("D29").Value = "'2.20"
$ws.Range("E29").Value = "  -1.09%  "
$ws.Range("D30").Value = "'166.30"
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("D31").Value = "'20.29"
$ws.Range("E31").Value = "  -3.64%  "
$ws.Range("D32").Value = "'0.0785"
$ws.Range("E32").Value = "  -2.25%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").Value = "'30.31"
$ws.Range("E33").Value = "  +3.09%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'5.32"
$ws.Range("E34").Value = "  -6.91%  "
$ws.Range("E35").Value = "  -1.80%  "
$ws.Range("D36").Value = "'0.108"
$ws.Range("E36").Value = "  -7.94%  "
$ws.Range("D37").Value = "'4.21"
$ws.Range("E37").Value = "  -4.88%  "
$ws.Range("D38").Value = "'0.0298"
$ws.Range("E38").Value = "  -2.02%  "
$ws.Range("D39").Value = "'12.87"
$ws.Range("E39").Value = "  -1.13%  "
$ws.Range("D40").Value = "'2.08"
$ws.Range("E40").Value = "  -3.57%  "
$ws.Range("D41").Value = "'5.57"
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("D42").Value = "'63.29"
$ws.Range("E42").Value = "  +0.56%  "
$ws.Range("D43").Value = "'0.194"
$ws.Range("E43").Value = "  -3.51%  "
$ws.Range("D44").Value = "'8.58"
$ws.Range("E44").Value = "  -2.68%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "'0.0991"
$ws.Range("E45").Value = "  -2.69%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'100.81"
$ws.Range("E46").Value = "  -4.23%  "
$ws.Range("D47").Value = "'1.11"
$ws.Range("E47").Value = "  -1.36%  "
$ws.Range("D48").Value = "'1.15"
$ws.Range("E48").Value = "  -1.88%  "
$ws.Range("D49").Value = "'2.29"
$ws.Range("E49").Value = "  -3.54%  "
$ws.Range("E50").Value = "  -1.56%  "
$ws.Range("D51").Value = "2.424.63"
$ws.Range("E51").Value = "  -0.41%  "
